# Daily attendance processing - 2025-12-26 09:34:22
#
# Normalises the "Recorded By" (column G) audit trail on the
# "Session Analysis Results" sheet: whichever account is listed FIRST in a
# multi-author cell gets rotated to the end of the list (i.e. the next
# editor in the chain moves to the front), except when the cell already
# starts with the shared "backup@backdoor.com" service account, which is
# left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByCol = 7   # column G - "Recorded By"
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    if ($parts[0] -eq "backup@backdoor.com") {
        continue
    }

    $rest = $parts[1..($parts.Count - 1)]
    $rotated = $rest + $parts[0]
    $newVal = $rotated -join ", "

    $cell.Value2 = $newVal
}
